$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price column D, Volume(1h) column E)
# Each cell is forced to Text format ("@") before assignment so the
# numeric-looking strings are stored as text, matching the source data.
$updates = @{
    "D2" = "332.85"
    "E2" = "1.01%"
    "D3" = "41.23"
    "E3" = "1.35%"
    "D4" = "5.694"
    "E4" = "-5.01%"
    "D5" = "0.08075"
    "E5" = "-1.51%"
    "D6" = "2.039"
    "E6" = "3.09%"
    "D7" = "8.746"
    "E7" = "-0.40%"
    "D8" = "4.540"
    "E8" = "-0.43%"
    "D9" = "2.966"
    "E9" = "-1.14%"
    "D10" = "0.9226"
    "E10" = "-3.03%"
    "D11" = "0.1266"
    "E11" = "-6.09%"
    "D12" = "0.1951"
    "E12" = "-2.43%"
    "D13" = "8.841"
    "E13" = "4.74%"
    "D14" = "0.09211"
    "E14" = "-0.45%"
    "D15" = "0.03682"
    "E15" = "5.40%"
    "D16" = "0.1052"
    "E16" = "8.53%"
    "D17" = "0.001299"
    "E17" = "-0.79%"
    "D18" = "0.006274"
    "E18" = "2.94%"
    "D19" = "3.373"
    "E19" = "0.32%"
    "D20" = "0.3481"
    "E20" = "-1.45%"
    "D21" = "0.1420"
    "E21" = "-1.65%"
    "E22" = "9.03%"
    "D23" = "0.04452"
    "E23" = "0.73%"
    "D24" = "0.001261"
    "E24" = "0.58%"
    "D25" = "0.004297"
    "E25" = "-3.68%"
    "D26" = "0.0001243"
    "E26" = "4.50%"
    "E39" = "15.81%"
    "D40" = "0.05502"
    "E40" = "3.24%"
    "D41" = "0.007789"
    "E41" = "4.01%"
    "D42" = "0.009932"
    "E42" = "11.01%"
    "D43" = "0.1419"
    "E43" = "-2.28%"
    "D44" = "0.002126"
    "E44" = "3.63%"
    "D45" = "0.01113"
    "E45" = "5.06%"
    "D46" = "0.00006826"
    "E46" = "0.95%"
    "D47" = "0.00000000752"
    "E47" = "0.64%"
    "D48" = "0.003022"
    "E48" = "-12.84%"
    "D49" = "0.002286"
    "E49" = "27.42%"
    "D50" = "0.00002106"
    "E50" = "0.64%"
    "D51" = "0.0002006"
    "E51" = "0.64%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}

Write-Output ("Updated " + $updates.Count + " cells")
